$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure values are written as text (matching original inlineStr cell type)
# rather than being auto-coerced to numbers by Excel, then restore default styling.
$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = '60.311.62'
$ws.Range("E2").Value = '  -4.13%  '
$ws.Range("D3").Value = '2.906.97'
$ws.Range("E3").Value = '  -3.48%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '527.94'
$ws.Range("E5").Value = '  -5.19%  '
$ws.Range("D6").Value = '142.22'
$ws.Range("E6").Value = '  -7.18%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '0.552'
$ws.Range("E8").Value = '  -1.97%  '
$ws.Range("D9").Value = '2.906.72'
$ws.Range("E9").Value = '  -3.65%  '
$ws.Range("E10").Value = '  -4.73%  '
$ws.Range("D11").Value = '5.86'
$ws.Range("E11").Value = '  -8.42%  '
$ws.Range("D12").Value = '0.353'
$ws.Range("E12").Value = '  -3.11%  '
$ws.Range("D13").Value = '3.414.29'
$ws.Range("E13").Value = '  -3.57%  '
$ws.Range("E14").Value = '  +1.24%  '
$ws.Range("D15").Value = '60.457.33'
$ws.Range("E15").Value = '  -4.03%  '
$ws.Range("D16").Value = '22.64'
$ws.Range("E16").Value = '  -5.48%  '
$ws.Range("D17").Value = '2.913.81'
$ws.Range("E17").Value = '  -3.37%  '
$ws.Range("D18").Value = '0.0000141'
$ws.Range("E18").Value = '  -5.96%  '
$ws.Range("D19").Value = '4.92'
$ws.Range("E19").Value = '  -3.57%  '
$ws.Range("D20").Value = '11.56'
$ws.Range("E20").Value = '  -2.57%  '
$ws.Range("D21").Value = '360.51'
$ws.Range("E21").Value = '  -8.82%  '
$ws.Range("E22").Value = '  -0.37%  '
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("D24").Value = '5.64'
$ws.Range("E24").Value = '  -1.89%  '
$ws.Range("D25").Value = '63.22'
$ws.Range("E25").Value = '  -2.91%  '
$ws.Range("D26").Value = '3.030.72'
$ws.Range("E27").Value = '  -3.42%  '
$ws.Range("E28").Value = '  -6.55%  '
$ws.Range("E29").Value = '  +0.42%  '
$ws.Range("D30").Value = '0.0₃0855'
$ws.Range("E30").Value = '  -12.03%  '
$ws.Range("D31").Value = '7.70'
$ws.Range("E31").Value = '  -11.01%  '
$ws.Range("E32").Value = '  -0.02%  '
$ws.Range("E33").Value = '  -4.66%  '
$ws.Range("D34").Value = '19.68'
$ws.Range("E34").Value = '  -3.76%  '
$ws.Range("D35").Value = '153.63'
$ws.Range("E35").Value = '  -4.02%  '
$ws.Range("E36").Value = '  -8.16%  '
$ws.Range("E37").Value = '  -8.05%  '
$ws.Range("D38").Value = '0.997'
$ws.Range("E38").Value = '  -8.91%  '
$ws.Range("D39").Value = '1.20'
$ws.Range("E39").Value = '  -7.82%  '
$ws.Range("D40").Value = '37.81'
$ws.Range("E40").Value = '  +0.65%  '
$ws.Range("D41").Value = '2.337.97'
$ws.Range("E41").Value = '  -6.78%  '
$ws.Range("E42").Value = '  -7.73%  '
$ws.Range("E43").Value = '  -6.27%  '
$ws.Range("E44").Value = '  -3.49%  '
$ws.Range("D45").Value = '20.86'
$ws.Range("E45").Value = '  -7.49%  '
$ws.Range("D46").Value = '0.0568'
$ws.Range("E46").Value = '  -4.90%  '
$ws.Range("E47").Value = '  -0.05%  '
$ws.Range("D48").Value = '4.85'
$ws.Range("E48").Value = '  -4.30%  '
$ws.Range("E49").Value = '  -1.11%  '
$ws.Range("E50").Value = '  -5.54%  '
$ws.Range("D51").Value = '0.0924'
$ws.Range("E51").Value = '  -2.34%  '

$rng.Style = "Normal"
